$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Alt")

# ---------------------------------------------------------------------------
# The "Alt" sheet re-tested its TEXTSPLIT/ANCHORARRAY based extraction
# approach against the B19 signature line ("John Doe Communications Manager
# john.doe@company.com") the same way it was already tested against B16,
# and that uncovered a bug in the alternate solution.  This mirrors the
# existing C16/C17/C18 block, one-for-one, for row 19/20/21.
# ---------------------------------------------------------------------------

# Row 19: split B19 into words (spills C19:G19)
$ws.Range("C19:G19").FormulaArray = "=TEXTSPLIT(B19, "" "")"

# Row 20: pull out the single "word" that looks like an email address
$ws.Range("C20").FormulaArray = "=TOCOL(IFS(1 - ISERR(FIND(""."",ANCHORARRAY( C19)) + FIND(""@"",ANCHORARRAY( C19))),ANCHORARRAY( C19)), 2)"

# Row 21: flag (0/1) which of the split words looks like an email address
$ws.Range("C21:G21").FormulaArray = "=1-ISERR(FIND(""."",ANCHORARRAY(C19))+FIND(""@"",ANCHORARRAY(C19)))"

# ---------------------------------------------------------------------------
# The consolidated MAP()/LAMBDA() extraction used to spill from B21:B29; it
# now starts two rows lower (B23:B31) to make room for the new C19:G21
# scratch-work above, with a blank spacer row (22) in between - mirroring
# the blank row 20 gap that already existed between rows 19 and 21.
# ---------------------------------------------------------------------------

# Clear out the old array (and its spill) before re-anchoring it lower down.
$ws.Range("B21:B29").Clear()

# B21:B29 used to carry an explicit style (index 3) instead of the column's
# own default (index 8); restore that same override on the new home range
# by cloning the formatting from a cell that already carries it.
$ws.Range("C16").Copy()
$ws.Range("B23:B31").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B23:B31").FormulaArray = "=TOCOL(MAP(B3:B19,LAMBDA(x,LET(y,TEXTSPLIT(x, "" ""),TOCOL(IFS(1 - ISERR(FIND(""."", y) + FIND(""@"", y)), y), 2)))),2)"

# ---------------------------------------------------------------------------
# Add a verification column (C) next to the extracted results, comparing
# each extracted address against the known-good address in column D (same
# pattern already used on the "EDA" sheet), plus a note for the rows where
# the alternate approach erroneously keeps a trailing period.
# ---------------------------------------------------------------------------

$ws.Range("C23").Formula = "=B23=D3"
$ws.Range("C24:C31").Formula = "=B24=D4"

$ws.Range("D28").Value = "Erroneously included period at the end"
$ws.Range("D29").Value = "Erroneously included period at the end"
$ws.Range("D30").Value = "Erroneously included period at the end"

# ---------------------------------------------------------------------------
# Update the view state to match where the author ended up looking.
# ---------------------------------------------------------------------------

$ws.Activate()
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("B36").Select()
